$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.866.90'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '1.977.49'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '245.06'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('E6').Value = '  +1.50%  '
$ws.Range('D7').Value = '60.78'
$ws.Range('E7').Value = '  +3.00%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +1.90%  '
$ws.Range('D10').Value = '0.0802'
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('E11').Value = '  +0.85%  '
$ws.Range('D12').Value = '14.66'
$ws.Range('E12').Value = '  +6.63%  '
$ws.Range('D13').Value = '0.845'
$ws.Range('E13').Value = '  +2.00%  '
$ws.Range('E14').Value = '  -1.56%  '
$ws.Range('D15').Value = '2.268.04'
$ws.Range('E15').Value = '  +0.50%  '
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('D17').Value = '1.974.85'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('D18').Value = '36.765.06'
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('D19').Value = '70.08'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('D20').Value = '0.0₃0861'
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').Value = '5.16'
$ws.Range('E21').Value = '  +1.45%  '
$ws.Range('D22').Value = '230.07'
$ws.Range('E22').Value = '  +0.39%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').Value = '2.49'
$ws.Range('E24').Value = '  +1.69%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('E26').Value = '  +1.97%  '
$ws.Range('D27').Value = '9.29'
$ws.Range('E27').Value = '  -0.28%  '
$ws.Range('D28').Value = '163.55'
$ws.Range('E28').Value = '  +1.89%  '
$ws.Range('D29').Value = '19.47'
$ws.Range('E29').Value = '  +0.32%  '
$ws.Range('D30').Value = '1.36'
$ws.Range('E30').Value = '  +19.64%  '
$ws.Range('E31').Value = '  +1.66%  '
$ws.Range('D32').Value = '4.85'
$ws.Range('E32').Value = '  +2.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0620'
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('E34').Value = '  +4.76%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('E37').Value = '  -1.88%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').Value = '5.47'
$ws.Range('E39').Value = '  -9.95%  '
$ws.Range('D40').Value = '0.0978'
$ws.Range('E40').Value = '  -2.35%  '
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('E42').Value = '  +0.85%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').Value = '16.21'
$ws.Range('E44').Value = '  +0.57%  '
$ws.Range('D45').Value = '1.368.10'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = '89.83'
$ws.Range('E46').Value = '  +2.21%  '
$ws.Range('E47').Value = '  -0.23%  '
$ws.Range('D48').Value = '7.24'
$ws.Range('E48').Value = '  +1.19%  '
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('D50').Value = '46.25'
$ws.Range('E50').Value = '  +5.58%  '
$ws.Range('D51').Value = '1.95'
$ws.Range('E51').Value = '  +9.27%  '
